$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "https://www.varoom.com/property/agriturismo-marano/EP-30538840"
$ws.Range("B3").Value = "https://www.varoom.com/property/residence-inn-indianapolis-fishers/BC-269249"
$ws.Range("B4").Value = "https://www.varoom.com/property/stone-soup-inn/EP-3852352"

$ws.Range("A5").Value = "www.varoom.com"
$ws.Range("B5").Value = "https://www.varoom.com/property/comfort-inn-suites-fishers-indianapolis/BC-183994"
$ws.Range("C5").Value = "Hybrid"
$ws.Range("D5").Value = "Property available in date range"
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = "The Property in the www.varoom.com is Available in the Specified date range"
